$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2117117117117117
$ws.Range("C2").Value = 0.527027027027027
$ws.Range("J2").Value = 0.009009009009009009
$ws.Range("P2").Value = 0.1486486486486487
$ws.Range("S2").Value = 0.1036036036036036
$ws.Range("B3").Value = 0.008130081300813009
$ws.Range("C3").Value = 0.02439024390243903
$ws.Range("J3").Value = 0.01626016260162602
$ws.Range("P3").Value = 0.6504065040650406
$ws.Range("S3").Value = 0.3008130081300813
$ws.Range("J4").Value = 0.02941176470588235
$ws.Range("P4").Value = 0.6764705882352942
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("B6").Value = 0.0273224043715847
$ws.Range("D6").Value = 0.03278688524590164
$ws.Range("F6").Value = 0.0546448087431694
$ws.Range("J6").Value = 0.2786885245901639
$ws.Range("O6").Value = 0.0273224043715847
$ws.Range("Q6").Value = 0.1366120218579235
$ws.Range("R6").Value = 0.07650273224043716
$ws.Range("S6").Value = 0.366120218579235
$ws.Range("B7").Value = 0.1219512195121951
$ws.Range("D7").Value = 0.01219512195121951
$ws.Range("F7").Value = 0.0426829268292683
$ws.Range("J7").Value = 0.1280487804878049
$ws.Range("O7").Value = 0.01829268292682927
$ws.Range("Q7").Value = 0.1402439024390244
$ws.Range("R7").Value = 0.07926829268292683
$ws.Range("S7").Value = 0.4573170731707317
$ws.Range("B8").Value = 0.08136482939632546
$ws.Range("D8").Value = 0.01312335958005249
$ws.Range("E8").Value = 0.002624671916010499
$ws.Range("F8").Value = 0.03937007874015748
$ws.Range("J8").Value = 0.1286089238845144
$ws.Range("O8").Value = 0.01049868766404199
$ws.Range("Q8").Value = 0.1653543307086614
$ws.Range("R8").Value = 0.09973753280839895
$ws.Range("S8").Value = 0.4593175853018373
$ws.Range("B9").Value = 0.0797872340425532
$ws.Range("D9").Value = 0.01595744680851064
$ws.Range("F9").Value = 0.06382978723404255
$ws.Range("J9").Value = 0.1170212765957447
$ws.Range("O9").Value = 0.02659574468085106
$ws.Range("Q9").Value = 0.1382978723404255
$ws.Range("R9").Value = 0.0851063829787234
$ws.Range("S9").Value = 0.4734042553191489
$ws.Range("B10").Value = 0.0890909090909091
$ws.Range("D10").Value = 0.01909090909090909
$ws.Range("E10").Value = 0.001818181818181818
$ws.Range("F10").Value = 0.07454545454545454
$ws.Range("J10").Value = 0.1218181818181818
$ws.Range("O10").Value = 0.01909090909090909
$ws.Range("Q10").Value = 0.1981818181818182
$ws.Range("R10").Value = 0.09454545454545454
$ws.Range("S10").Value = 0.3818181818181818
$ws.Range("G11").Value = 0.1245283018867925
$ws.Range("J11").Value = 0.09811320754716982
$ws.Range("K11").Value = 0.1849056603773585
$ws.Range("L11").Value = 0.569811320754717
$ws.Range("S11").Value = 0.02264150943396226
$ws.Range("G12").Value = 0.732484076433121
$ws.Range("J12").Value = 0.2038216560509554
$ws.Range("K12").Value = 0.006369426751592357
$ws.Range("L12").Value = 0.02547770700636943
$ws.Range("S12").Value = 0.03184713375796178
$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.3714285714285714
$ws.Range("S13").Value = 0.05714285714285714
$ws.Range("F15").Value = 0.01136363636363636
$ws.Range("H15").Value = 0.1363636363636364
$ws.Range("I15").Value = 0.07386363636363637
$ws.Range("J15").Value = 0.3238636363636364
$ws.Range("K15").Value = 0.05681818181818182
$ws.Range("M15").Value = 0.01704545454545454
$ws.Range("O15").Value = 0.06818181818181818
$ws.Range("S15").Value = 0.3125
$ws.Range("F16").Value = 0.02238805970149254
$ws.Range("H16").Value = 0.1492537313432836
$ws.Range("I16").Value = 0.08955223880597014
$ws.Range("J16").Value = 0.417910447761194
$ws.Range("K16").Value = 0.1194029850746269
$ws.Range("M16").Value = 0.007462686567164179
$ws.Range("O16").Value = 0.03731343283582089
$ws.Range("S16").Value = 0.1567164179104478
$ws.Range("F17").Value = 0.01412429378531073
$ws.Range("H17").Value = 0.1355932203389831
$ws.Range("I17").Value = 0.09322033898305085
$ws.Range("J17").Value = 0.4576271186440678
$ws.Range("K17").Value = 0.09322033898305085
$ws.Range("M17").Value = 0.01129943502824859
$ws.Range("N17").Value = 0.002824858757062147
$ws.Range("O17").Value = 0.06214689265536723
$ws.Range("S17").Value = 0.1299435028248588
$ws.Range("F18").Value = 0.01630434782608696
$ws.Range("H18").Value = 0.1630434782608696
$ws.Range("I18").Value = 0.09239130434782608
$ws.Range("J18").Value = 0.3695652173913043
$ws.Range("K18").Value = 0.1304347826086956
$ws.Range("M18").Value = 0.0108695652173913
$ws.Range("O18").Value = 0.07065217391304347
$ws.Range("S18").Value = 0.1467391304347826
$ws.Range("F19").Value = 0.01526717557251908
$ws.Range("H19").Value = 0.2222222222222222
$ws.Range("I19").Value = 0.09923664122137404
$ws.Range("J19").Value = 0.356234096692112
$ws.Range("K19").Value = 0.1085665818490246
$ws.Range("M19").Value = 0.02205258693808312
$ws.Range("N19").Value = 0.001696352841391009
$ws.Range("O19").Value = 0.05767599660729432
$ws.Range("S19").Value = 0.1170483460559796
